$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix CTRL-001's description: the old wording described a fully automated
# disable after 90 days (a precise WHEN), but the control as implemented is
# really a daily review/manual-disable step, so the vague-timing text is
# replaced to correctly flag it as "missing" a precise WHEN under the WHEN
# scoring rubric.
$ws.Range("B2").Value = "The Accounting Manager reviews monthly bank reconciliations prepared by the Senior Accountant to ensure completeness and accuracy. Reconciliations are completed by the 5th business day. Unresolved items over `$1,000 are escalated to the Controller."

# Append a new control row (CTRL-011) describing a vague "monthly" review
# with no explicit calendar trigger.
$ws.Range("A12").Value = "CTRL-011"
$ws.Range("B12").Value = "On a monthly basis, the Finance Accounts Receivable Manager reviews delinquent balances over 120 days old to verify compliance with the Company's write-off policy and regulatory requirements. The manager validates that all qualifying accounts have been properly identified for write-off by comparing the aged receivables report against established thresholds. Any accounts incorrectly processed are documented, and discrepancies exceeding `$10,000 are escalated to the Finance Accounts Receivable Director for resolution."
$ws.Range("C12").Value = "Monthly"
$ws.Range("D12").Value = "Semi-Automated"
$ws.Range("E12").Value = "Maria Chen"

$ws.Range("D12").Select()
